$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of search-result data appended below the existing rows (19 -> 22).
# Columns: A=Date(serial), B=Method(string "Named"), C..M numeric metrics.
$newRows = @(
    @(42602.58184027778,  "Named", 10245, 8234, 510, 66, 48, 57, 41, 0, 0,   0, 0),
    @(42602.980011574073, "Named", 8861,  5187, 299, 47, 28, 62, 37, 0, 0,   0, 0),
    @(42603.694513888891, "Named", 5238,  1974, 132, 20, 10, 66, 33, 1, 0, 100, 0)
)

$startRow = 20
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
